# Auto update Excel log 2026-02-04 14:20:19
# Appends newly-logged sensor readings to the PIR, Humidity and Temperature
# sheets of the SeniorConnect master log workbook. Each new record uses the
# same 6-column (Date, Timestamp, Hour, Location, Value, Status) shape as
# the pre-existing rows. Values are written as literal text (not an
# auto-converted date/time/number) by forcing a "@" text format before the
# values are assigned, then restoring the default "Normal" style so the
# appended cells render identically to the rest of the log.

$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 217-230 ---
$ws = $wb.Worksheets.Item("PIR")
$newRows = @(
    @{ Row = 217; A = '2026-02-04'; B = '14:19:16'; C = '14:00'; D = 'Bathroom'; E = 'No Motion'; F = 'Inactive' },
    @{ Row = 218; A = '2026-02-04'; B = '14:19:17'; C = '14:00'; D = 'Bathroom'; E = 'Motion Detected'; F = 'Active' },
    @{ Row = 219; A = '2026-02-04'; B = '14:19:19'; C = '14:00'; D = 'Bathroom'; E = 'No Motion'; F = 'Inactive' },
    @{ Row = 220; A = '2026-02-04'; B = '14:19:24'; C = '14:00'; D = 'Bathroom'; E = 'No Motion'; F = 'Inactive' },
    @{ Row = 221; A = '2026-02-04'; B = '14:19:25'; C = '14:00'; D = 'Bathroom'; E = 'Motion Detected'; F = 'Active' },
    @{ Row = 222; A = '2026-02-04'; B = '14:19:34'; C = '14:00'; D = 'Bathroom'; E = 'No Motion'; F = 'Inactive' },
    @{ Row = 223; A = '2026-02-04'; B = '14:19:40'; C = '14:00'; D = 'Bathroom'; E = 'No Motion'; F = 'Inactive' },
    @{ Row = 224; A = '2026-02-04'; B = '14:19:44'; C = '14:00'; D = 'Bathroom'; E = 'No Motion'; F = 'Inactive' },
    @{ Row = 225; A = '2026-02-04'; B = '14:19:45'; C = '14:00'; D = 'Bathroom'; E = 'Motion Detected'; F = 'Active' },
    @{ Row = 226; A = '2026-02-04'; B = '14:19:52'; C = '14:00'; D = 'Bathroom'; E = 'No Motion'; F = 'Inactive' },
    @{ Row = 227; A = '2026-02-04'; B = '14:19:57'; C = '14:00'; D = 'Bathroom'; E = 'No Motion'; F = 'Inactive' },
    @{ Row = 228; A = '2026-02-04'; B = '14:20:02'; C = '14:00'; D = 'Bathroom'; E = 'No Motion'; F = 'Inactive' },
    @{ Row = 229; A = '2026-02-04'; B = '14:20:07'; C = '14:00'; D = 'Bathroom'; E = 'No Motion'; F = 'Inactive' },
    @{ Row = 230; A = '2026-02-04'; B = '14:20:09'; C = '14:00'; D = 'Bathroom'; E = 'Motion Detected'; F = 'Active' }
)
foreach ($item in $newRows) {
    $r = $ws.Range("A" + $item.Row + ":F" + $item.Row)
    $r.NumberFormat = "@"
    $ws.Range("A" + $item.Row).Value = $item.A
    $ws.Range("B" + $item.Row).Value = $item.B
    $ws.Range("C" + $item.Row).Value = $item.C
    $ws.Range("D" + $item.Row).Value = $item.D
    $ws.Range("E" + $item.Row).Value = $item.E
    $ws.Range("F" + $item.Row).Value = $item.F
    $r.Style = "Normal"
}

# --- Humidity sheet: append rows 182-194 ---
$ws = $wb.Worksheets.Item("Humidity")
$newRows = @(
    @{ Row = 182; A = '2026-02-04'; B = '14:19:15'; C = '14:00'; D = 'Bathroom'; E = '79.1%'; F = 'Active' },
    @{ Row = 183; A = '2026-02-04'; B = '14:19:18'; C = '14:00'; D = 'Bathroom'; E = '77.7%'; F = 'Active' },
    @{ Row = 184; A = '2026-02-04'; B = '14:19:23'; C = '14:00'; D = 'Bathroom'; E = '78.2%'; F = 'Active' },
    @{ Row = 185; A = '2026-02-04'; B = '14:19:28'; C = '14:00'; D = 'Bathroom'; E = '79.1%'; F = 'Active' },
    @{ Row = 186; A = '2026-02-04'; B = '14:19:33'; C = '14:00'; D = 'Bathroom'; E = '78.3%'; F = 'Active' },
    @{ Row = 187; A = '2026-02-04'; B = '14:19:38'; C = '14:00'; D = 'Bathroom'; E = '79.2%'; F = 'Active' },
    @{ Row = 188; A = '2026-02-04'; B = '14:19:43'; C = '14:00'; D = 'Bathroom'; E = '78.2%'; F = 'Active' },
    @{ Row = 189; A = '2026-02-04'; B = '14:19:48'; C = '14:00'; D = 'Bathroom'; E = '79.0%'; F = 'Active' },
    @{ Row = 190; A = '2026-02-04'; B = '14:19:53'; C = '14:00'; D = 'Bathroom'; E = '79.6%'; F = 'Active' },
    @{ Row = 191; A = '2026-02-04'; B = '14:19:58'; C = '14:00'; D = 'Bathroom'; E = '80.3%'; F = 'Active' },
    @{ Row = 192; A = '2026-02-04'; B = '14:20:03'; C = '14:00'; D = 'Bathroom'; E = '78.4%'; F = 'Active' },
    @{ Row = 193; A = '2026-02-04'; B = '14:20:08'; C = '14:00'; D = 'Bathroom'; E = '79.3%'; F = 'Active' },
    @{ Row = 194; A = '2026-02-04'; B = '14:20:14'; C = '14:00'; D = 'Bathroom'; E = '78.1%'; F = 'Active' }
)
foreach ($item in $newRows) {
    $r = $ws.Range("A" + $item.Row + ":F" + $item.Row)
    $r.NumberFormat = "@"
    $ws.Range("A" + $item.Row).Value = $item.A
    $ws.Range("B" + $item.Row).Value = $item.B
    $ws.Range("C" + $item.Row).Value = $item.C
    $ws.Range("D" + $item.Row).Value = $item.D
    $ws.Range("E" + $item.Row).Value = $item.E
    $ws.Range("F" + $item.Row).Value = $item.F
    $r.Style = "Normal"
}

# --- Temperature sheet: append rows 182-194 ---
$ws = $wb.Worksheets.Item("Temperature")
$newRows = @(
    @{ Row = 182; A = '2026-02-04'; B = '14:19:15'; C = '14:00'; D = 'Bathroom'; E = '24.4C'; F = 'Active' },
    @{ Row = 183; A = '2026-02-04'; B = '14:19:19'; C = '14:00'; D = 'Bathroom'; E = '24.5C'; F = 'Active' },
    @{ Row = 184; A = '2026-02-04'; B = '14:19:24'; C = '14:00'; D = 'Bathroom'; E = '24.4C'; F = 'Active' },
    @{ Row = 185; A = '2026-02-04'; B = '14:19:29'; C = '14:00'; D = 'Bathroom'; E = '24.5C'; F = 'Active' },
    @{ Row = 186; A = '2026-02-04'; B = '14:19:34'; C = '14:00'; D = 'Bathroom'; E = '24.5C'; F = 'Active' },
    @{ Row = 187; A = '2026-02-04'; B = '14:19:39'; C = '14:00'; D = 'Bathroom'; E = '24.4C'; F = 'Active' },
    @{ Row = 188; A = '2026-02-04'; B = '14:19:44'; C = '14:00'; D = 'Bathroom'; E = '24.5C'; F = 'Active' },
    @{ Row = 189; A = '2026-02-04'; B = '14:19:49'; C = '14:00'; D = 'Bathroom'; E = '24.4C'; F = 'Active' },
    @{ Row = 190; A = '2026-02-04'; B = '14:19:54'; C = '14:00'; D = 'Bathroom'; E = '24.5C'; F = 'Active' },
    @{ Row = 191; A = '2026-02-04'; B = '14:19:59'; C = '14:00'; D = 'Bathroom'; E = '24.4C'; F = 'Active' },
    @{ Row = 192; A = '2026-02-04'; B = '14:20:04'; C = '14:00'; D = 'Bathroom'; E = '24.5C'; F = 'Active' },
    @{ Row = 193; A = '2026-02-04'; B = '14:20:09'; C = '14:00'; D = 'Bathroom'; E = '24.5C'; F = 'Active' },
    @{ Row = 194; A = '2026-02-04'; B = '14:20:14'; C = '14:00'; D = 'Bathroom'; E = '24.5C'; F = 'Active' }
)
foreach ($item in $newRows) {
    $r = $ws.Range("A" + $item.Row + ":F" + $item.Row)
    $r.NumberFormat = "@"
    $ws.Range("A" + $item.Row).Value = $item.A
    $ws.Range("B" + $item.Row).Value = $item.B
    $ws.Range("C" + $item.Row).Value = $item.C
    $ws.Range("D" + $item.Row).Value = $item.D
    $ws.Range("E" + $item.Row).Value = $item.E
    $ws.Range("F" + $item.Row).Value = $item.F
    $r.Style = "Normal"
}
